$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-12-19 Friday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-12-20 Saturday", 2)
$d.Content.Find.Execute("661÷2=330, 1", $true, $false, $false, $false, $false, $true, 1, $false, "792÷6=132, 0", 2)
$d.Content.Find.Execute("250÷6=41, 4", $true, $false, $false, $false, $false, $true, 1, $false, "744÷6=124, 0", 2)
$d.Content.Find.Execute("408÷6=68, 0", $true, $false, $false, $false, $false, $true, 1, $false, "264÷8=33, 0", 2)
$d.Content.Find.Execute("476÷2=238, 0", $true, $false, $false, $false, $false, $true, 1, $false, "744÷6=124, 0", 2)
$d.Content.Find.Execute("155÷7=22, 1", $true, $false, $false, $false, $false, $true, 1, $false, "260÷9=28, 8", 2)
$d.Content.Find.Execute("786÷9=87, 3", $true, $false, $false, $false, $false, $true, 1, $false, "789÷4=197, 1", 2)
$d.Content.Find.Execute("300÷4=75, 0", $true, $false, $false, $false, $false, $true, 1, $false, "803÷8=100, 3", 2)
$d.Content.Find.Execute("139÷4=34, 3", $true, $false, $false, $false, $false, $true, 1, $false, "518÷4=129, 2", 2)
$d.Content.Find.Execute("467÷5=93, 2", $true, $false, $false, $false, $false, $true, 1, $false, "379÷3=126, 1", 2)
$d.Content.Find.Execute("285÷4=71, 1", $true, $false, $false, $false, $false, $true, 1, $false, "446÷5=89, 1", 2)
$d.Content.Find.Execute("306÷9=34, 0", $true, $false, $false, $false, $false, $true, 1, $false, "477÷3=159, 0", 2)
$d.Content.Find.Execute("188÷6=31, 2", $true, $false, $false, $false, $false, $true, 1, $false, "537÷2=268, 1", 2)
$d.Content.Find.Execute("468÷9=52, 0", $true, $false, $false, $false, $false, $true, 1, $false, "787÷8=98, 3", 2)
$d.Content.Find.Execute("172÷7=24, 4", $true, $false, $false, $false, $false, $true, 1, $false, "561÷4=140, 1", 2)
$d.Content.Find.Execute("912÷7=130, 2", $true, $false, $false, $false, $false, $true, 1, $false, "103÷2=51, 1", 2)
$d.Content.Find.Execute("374÷5=74, 4", $true, $false, $false, $false, $false, $true, 1, $false, "662÷9=73, 5", 2)
$d.Content.Find.Execute("821÷6=136, 5", $true, $false, $false, $false, $false, $true, 1, $false, "937÷3=312, 1", 2)
$d.Content.Find.Execute("663÷5=132, 3", $true, $false, $false, $false, $false, $true, 1, $false, "274÷6=45, 4", 2)
$d.Content.Find.Execute("238÷9=26, 4", $true, $false, $false, $false, $false, $true, 1, $false, "723÷2=361, 1", 2)
$d.Content.Find.Execute("691÷8=86, 3", $true, $false, $false, $false, $false, $true, 1, $false, "560÷2=280, 0", 2)
$d.Content.Find.Execute("279÷7=39, 6", $true, $false, $false, $false, $false, $true, 1, $false, "514÷7=73, 3", 2)
$d.Content.Find.Execute("478÷4=119, 2", $true, $false, $false, $false, $false, $true, 1, $false, "956÷7=136, 4", 2)
$d.Content.Find.Execute("373÷7=53, 2", $true, $false, $false, $false, $false, $true, 1, $false, "592÷4=148, 0", 2)
$d.Content.Find.Execute("434÷7=62, 0", $true, $false, $false, $false, $false, $true, 1, $false, "681÷8=85, 1", 2)
$d.Content.Find.Execute("583÷4=145, 3", $true, $false, $false, $false, $false, $true, 1, $false, "837÷7=119, 4", 2)
